$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, pushing existing rows 66-299 down to 67-300.
$ws.Rows.Item(66).EntireRow.Insert()

# Populate the newly inserted row 66 with a new data record (same shape as the
# other rows in this table), using the same constant values shared by every
# row plus the new date / volume figures from this edit.
$ws.Cells.Item(66, 1).Value = 3
$ws.Cells.Item(66, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(66, 3).Value = "Coquimbo"
$ws.Cells.Item(66, 4).Value = "2022-04-08"
$ws.Cells.Item(66, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(66, 5).Value = 5
$ws.Cells.Item(66, 6).Value = 100112039
$ws.Cells.Item(66, 7).Value = "Ciboulette"
$ws.Cells.Item(66, 8).Value = "Sin especificar"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 120
$ws.Cells.Item(66, 11).Value = 1500
$ws.Cells.Item(66, 12).Value = 1500
$ws.Cells.Item(66, 13).Value = 1500
$ws.Cells.Item(66, 14).Value = "`$/docena de atados"
$ws.Cells.Item(66, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(66, 16).Value = 500
$ws.Cells.Item(66, 17).Value = 3
$ws.Cells.Item(66, 18).Value = "Hortaliza"
